# Apply "Changed Spell Dictionary" edits to spellDictionary.xlsx
# - Renames the basic physical-attack weapons (spear -> lance, falchion -> hammer)
# - Renames the elemental "adds X element to spell" entries to their new
#   lore names (fire -> agni, ice -> cryo, volt -> veld)
# - Rewrites the spellbook descriptions for those six entries
# - Leaves the selection on D28 to mirror the saved view state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename spell/weapon names (column A)
$ws.Range("A3").Value = "lance"
$ws.Range("A4").Value = "hammer"
$ws.Range("A25").Value = "agni"
$ws.Range("A26").Value = "cryo"
$ws.Range("A27").Value = "veld"

# Rewrite spellbook descriptions (column C)
$ws.Range("C2").Value = "Basic Phys ATK"
$ws.Range("C3").Value = "Medium Phys ATK"
$ws.Range("C4").Value = "Heavy Phys ATK"
$ws.Range("C25").Value = "Override ATK type to Fire"
$ws.Range("C26").Value = "Override ATK type to Ice"
$ws.Range("C27").Value = "Override ATK type to Volt"

# Match the final selection recorded in the saved workbook
$ws.Range("D28").Select()
